# Update the Handoff/Handback datetime stamps for the ab33e234 (...41675efc...)
# rows on the zh-cn and de-de language sheets, as part of regenerating the
# handback status report.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 3 holds the ab33e234.../zh-cn.xlf entry
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-28 04:08:10"
$wsZhCn.Range("G3").Value = "2016-01-28 04:08:50"

# de-de sheet: row 3 holds the ab33e234.../de-de.xlf entry
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-28 04:08:20"
$wsDeDe.Range("G3").Value = "2016-01-28 04:09:07"
